$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3.13
    3  = 3.35
    4  = 3.22
    5  = 3.31
    6  = 3.2
    7  = 3.32
    8  = 3.2
    9  = 3.07
    10 = 3.29
    11 = 3.19
    12 = 3.19
    13 = 3.08
    14 = 3.08
    15 = 3.21
    16 = 3.23
    17 = 3.11
    18 = 3.05
    19 = 3.05
    20 = 1.9
    21 = 2.33
    22 = 2.34
    23 = 2.3
    24 = 2.3
    25 = 2.36
    26 = 7.01
    27 = 3.62
    28 = 3.75
    29 = 2.7
    30 = 3.61
    31 = 8.51
    32 = 2.54
    33 = 2.57
    34 = 3.04
    35 = 3.67
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
